$wb = $excel.ActiveWorkbook

# --- Summary sheet ---
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B6").Value = 63
$wsSummary.Range("B9").Value = 39.68

# --- Strategy Status sheet ---
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("D4").Value = 63
$wsStatus.Range("G4").Value = 39.68

# --- New trade row data (Trade #63) ---
$newRow = @(63, "2026-02-17", "08:48:33", "MarketMaking", "UP", 0.8, 0.8, "CLOSED", 0, 0, 99.65000000000001, 0, 0, 0.6, "Normal spread capture: 19600 bps", "early_exit", 0.14)

# --- All Trades sheet ---
$wsAllTrades = $wb.Worksheets.Item("All Trades")
for ($i = 0; $i -lt $newRow.Length; $i++) {
    $wsAllTrades.Cells.Item(64, $i + 1).Value = $newRow[$i]
}
# Prevent the "2026-02-17" date-literal from being auto-converted into a
# date serial number by Excel's type inference: format as Text first so it
# is stored verbatim, then strip the now-unneeded format override so the
# cell is left with the default (unstyled) appearance, matching the rest
# of the sheet.
$wsAllTrades.Cells.Item(64, 2).NumberFormat = "@"
$wsAllTrades.Cells.Item(64, 2).Value = $newRow[1]
$wsAllTrades.Cells.Item(64, 2).ClearFormats()

# --- MarketMaking sheet ---
$wsMarketMaking = $wb.Worksheets.Item("MarketMaking")
for ($i = 0; $i -lt $newRow.Length; $i++) {
    $wsMarketMaking.Cells.Item(64, $i + 1).Value = $newRow[$i]
}
$wsMarketMaking.Cells.Item(64, 2).NumberFormat = "@"
$wsMarketMaking.Cells.Item(64, 2).Value = $newRow[1]
$wsMarketMaking.Cells.Item(64, 2).ClearFormats()
